$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly driver report update for 2025-04-21
$ws.Range("C3").Value = 2194
$ws.Range("D3").Value = 87.7

$ws.Range("D4").Value = 95.5

$ws.Range("C5").Value = 2259

$ws.Range("C6").Value = 294

$ws.Range("C7").Value = 4850
